$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers for columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy header style (bold/border/centered) from existing header cell H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the data for columns I (I0) and J (IF), rows 2-42
$data = @(
  @(8,9),
  @(9,9),
  @(8,8),
  @(9,9),
  @(9,9),
  @(7,7),
  @(6,6),
  @(7,7),
  @(7,7),
  @(8,8),
  @(8,8),
  @(8,8),
  @(8,8),
  @(6,6),
  @(8,8),
  @(6,6),
  @(11,11),
  @(7,7),
  @(6,6),
  @(6,7),
  @(8,8),
  @(8,8),
  @(6,6),
  @(7,7),
  @(8,8),
  @(7,7),
  @(11,11),
  @(8,8),
  @(6,6),
  @(7,7),
  @(8,8),
  @(7,7),
  @(7,7),
  @(8,8),
  @(7,7),
  @(7,7),
  @(7,7),
  @(7,7),
  @(8,8),
  @(8,8),
  @(6,6)
)

for ($n = 0; $n -lt $data.Length; $n++) {
  $r = $n + 2
  $ws.Cells.Item($r, 9).Value = $data[$n][0]
  $ws.Cells.Item($r, 10).Value = $data[$n][1]
}

Write-Output "Added I0 and IF columns"
